$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 144
$ws.Range("J2").Value = 495
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 136
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 95
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("R2").Value = 6
$ws.Range("T2").Value = 84
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 855
$ws.Range("X2").Value = 834
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 9
